# Init on emissions data handling
# Adds an "Atkvitetsdata" / "Emissionsfaktor" header pair to the "Scope 3"
# sheet, bolds + widens the new header row, and leaves the UI focused on
# the "Scope 3" sheet (matching the selection/activeTab changes in the diff).

$wb = $excel.ActiveWorkbook

# --- Scope 3 sheet: extend the header row from A1:B1 to A1:D1 -------------
$ws3 = $wb.Worksheets.Item("Scope 3")

$ws3.Range("C1").Value = "Atkvitetsdata"
$ws3.Range("D1").Value = "Emissionsfaktor"

# Header row is bold, like the existing A1/B1 header cells on the other
# scope sheets.
$ws3.Range("A1:D1").Font.Bold = $true

# Resize the columns to fit the new headers (values chosen so the stored,
# pixel-quantized column width lands as close as possible to the target).
$ws3.Columns.Item(1).ColumnWidth = 17.7214
$ws3.Columns.Item(2).ColumnWidth = 26.9440
$ws3.Columns.Item(3).ColumnWidth = 11.3854
$ws3.Columns.Item(4).ColumnWidth = 13.4987

# --- Selection / active-sheet bookkeeping ---------------------------------
# Leave "Scope 2" with a plain selection on F1 (no longer the active tab).
$ws2 = $wb.Worksheets.Item("Scope 2")
$ws2.Activate()
$ws2.Range("F1").Select() | Out-Null

# "Scope 3" becomes the active/selected sheet, with D6 selected.
$ws3.Activate()
$ws3.Range("D6").Select() | Out-Null
